# Insert a new data row right above the former row 272 (pushing the
# existing rows 272:367 down to 273:368) and populate it with a new
# weekly price observation for Mango at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 272:367 down by inserting a brand new row 272. Excel copies
# the formatting (e.g. the date number format on column D) from the row
# above, matching the style used throughout the rest of the table.
$ws.Rows.Item(272).Insert()

$ws.Cells.Item(272, 1).Value2 = 9
$ws.Cells.Item(272, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(272, 3).Value2 = "Metropolitana"
$ws.Cells.Item(272, 4).Value2 = 44627
$ws.Cells.Item(272, 5).Value2 = 13
$ws.Cells.Item(272, 6).Value2 = "Fruta"
$ws.Cells.Item(272, 7).Value2 = 100108
$ws.Cells.Item(272, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(272, 9).Value2 = 100108002
$ws.Cells.Item(272, 10).Value2 = "Mango"
$ws.Cells.Item(272, 11).Value2 = "Sin especificar"
$ws.Cells.Item(272, 12).Value2 = "Primera"
$ws.Cells.Item(272, 13).Value2 = 630
$ws.Cells.Item(272, 14).Value2 = 6000
$ws.Cells.Item(272, 15).Value2 = 6500
$ws.Cells.Item(272, 16).Value2 = 6302
$ws.Cells.Item(272, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(272, 18).Value2 = "Perú"
$ws.Cells.Item(272, 19).Value2 = 1576
$ws.Cells.Item(272, 20).Value2 = 4
